# Apply data refresh updates to the "南宁-漫展信息" workbook
# (mirrors the output regenerated at commit 456a3b4 for gh-pages)

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F2").Value = 5526
$wsExpo.Range("I2").Value = "//i0.hdslb.com/bfs/openplatform/202407/ilMDMXk11721378263909.jpeg"

$wsExpo.Range("F3").Value = 608
$wsExpo.Range("F4").Value = 12270
$wsExpo.Range("F7").Value = 186
$wsExpo.Range("F8").Value = 346
$wsExpo.Range("F9").Value = 1131
$wsExpo.Range("F10").Value = 107

# --- Sheet "全部类型" (all types) ---
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F3").Value = 5526
$wsAll.Range("I3").Value = "//i0.hdslb.com/bfs/openplatform/202407/ilMDMXk11721378263909.jpeg"

$wsAll.Range("F4").Value = 608
$wsAll.Range("F6").Value = 12270
$wsAll.Range("F9").Value = 186
$wsAll.Range("F12").Value = 346
$wsAll.Range("F13").Value = 1131
$wsAll.Range("F15").Value = 107
